# Applies two kinds of changes found in the commit diff:
#  1. Three tables (on slides 14, 15, 16) switch from the default
#     "Light Style 1" table style to the "Light Style 1 - Accent 2" style.
#  2. The deck's theme colour palette is swapped from the custom
#     "Integral" / "Red Violet" palette to the stock Office 2016-19
#     "Office Theme" palette (this is what ships as ppt/theme/theme1.xml,
#     the theme actually driving the slide master / all slides).

$p = $ppt.ActivePresentation

# --- 1. Re-colour the slide master theme (ppt/theme/theme1.xml) -----------
# Slide.ColorScheme exposes the 12 theme colours (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) in that order, regardless of which slide it
# is read from (there is a single slide master/theme in this deck).
$s1 = $p.Slides.Item(1)
$cs = $s1.ColorScheme

# RGB() packs as 0xBBGGRR, matching the legacy ColorScheme.Colors(i).RGB
# setter used by PowerPoint automation.
function Set-ThemeColor($scheme, $index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    $bgr = ($b * 65536) + ($g * 256) + $r
    $scheme.Colors($index).RGB = $bgr
}

Set-ThemeColor $cs 1  "000000"   # dk1
Set-ThemeColor $cs 2  "FFFFFF"   # lt1
Set-ThemeColor $cs 3  "44546A"   # dk2
Set-ThemeColor $cs 4  "E7E6E6"   # lt2
Set-ThemeColor $cs 5  "5B9BD5"   # accent1
Set-ThemeColor $cs 6  "ED7D31"   # accent2
Set-ThemeColor $cs 7  "A5A5A5"   # accent3
Set-ThemeColor $cs 8  "FFC000"   # accent4
Set-ThemeColor $cs 9  "4472C4"   # accent5
Set-ThemeColor $cs 10 "70AD47"   # accent6
Set-ThemeColor $cs 11 "0563C1"   # hlink
Set-ThemeColor $cs 12 "954F72"   # folHlink

# --- 2. Re-style the three tables on slides 14-16 --------------------------
$tableSlideIndexes = 14, 15, 16
foreach ($idx in $tableSlideIndexes) {
    $slide = $p.Slides.Item($idx)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle("{4D9FE028-4175-4768-9AF4-B05958465913}")
        }
    }
}
